$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B17:G17").Copy()
$ws.Range("B18:G18").PasteSpecial(-4122)
$ws.Range("B18").Value = 42181
$ws.Range("C18").Value = "Jovanny Zepeda"
$ws.Range("E18").Value = 42547
$ws.Range("G18").Value = "Realizada"
